$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 173 (shifts existing rows 173-183 down to 174-184),
# matching the new weekly "Femacal de La Calera - Achicoria" price record.
$ws.Rows.Item(173).Insert()

$ws.Cells.Item(173, 1).Value  = 3
$ws.Cells.Item(173, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(173, 3).Value  = "Coquimbo"
$ws.Cells.Item(173, 4).Value  = 44585
$ws.Cells.Item(173, 5).Value  = 5
$ws.Cells.Item(173, 6).Value  = 100112010
$ws.Cells.Item(173, 7).Value  = "Achicoria"
$ws.Cells.Item(173, 8).Value  = "Sin especificar"
$ws.Cells.Item(173, 9).Value  = "Primera"
$ws.Cells.Item(173, 10).Value = 95
$ws.Cells.Item(173, 11).Value = 5500
$ws.Cells.Item(173, 12).Value = 6000
$ws.Cells.Item(173, 13).Value = 5789
$ws.Cells.Item(173, 14).Value = "$/caja 16 unidades"
$ws.Cells.Item(173, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(173, 16).Value = 362
$ws.Cells.Item(173, 17).Value = 16
$ws.Cells.Item(173, 18).Value = "Hortaliza"
